# Append new rows (209-223) of NSAA measurement data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("NSAA", "position", "dhc", 60),
    @("NSAA", "position", "overall", 60),
    @("NSAA", "position", "acts", 60),
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "sensorMagneticField", "overall", 60),
    @("NSAA", "sensorMagneticField", "acts", 60),
    @("NSAA", "jointAngle", "dhc", 60),
    @("NSAA", "jointAngle", "overall", 60),
    @("NSAA", "jointAngle", "acts", 60),
    @("NSAA", "jointAngleXZY", "dhc", 60),
    @("NSAA", "jointAngleXZY", "overall", 60),
    @("NSAA", "jointAngleXZY", "acts", 60),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "AD", "acts", 10)
)

$startRow = 209
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
